# Applies the "New crime data collected" weekly update to the 116th Precinct
# CompStat workbook.
#
# Summary of changes (per commit diff):
#  - Report header: Volume/Number bumped from "...29" to "...30"
#  - Report header: week-covering dates shifted forward one week
#  - Crime-table rows 15-30 (columns C:L) updated with the new week's figures,
#    including a couple of cells that flip between numeric 0/"N/A" style
#    (displayed as "0" / "***.*") and real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: Volume 32  Number 29 -> 30, and the covered-week date range
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  30"
$ws.Range("C9").Value = "Report Covering the Week  7/21/2025  Through  7/27/2025"

# ---------------------------------------------------------------------------
# Helper: write a "blank"/placeholder cell such as the ones displayed as
# "0" or "***.*" (these are shared-text placeholders, not real numbers).
# Style is copied from the neighboring placeholder cell in column C so the
# numeric formatting (right aligned, General number format) matches the
# other placeholder cells on the sheet.
# ---------------------------------------------------------------------------
function Set-Placeholder($cellRef, $text) {
    $target = $ws.Range($cellRef)
    $target.Value = "'" + $text
    $ws.Range("C15").Copy() | Out-Null
    $target.PasteSpecial(-4122) | Out-Null
}

# Row 15 (Murder): D/E flip from real numbers to the "0" / "***.*" placeholders
Set-Placeholder "D15" "0"
Set-Placeholder "E15" "***.*"

# Row 19 (Burglary): C flips from a real number to the "0" placeholder
Set-Placeholder "C19" "0"

# Row 27 (UCR Rape*): D/E flip from real numbers to the "0" / "***.*" placeholders
Set-Placeholder "D27" "0"
Set-Placeholder "E27" "***.*"

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Crime statistics table: updated weekly figures
# ---------------------------------------------------------------------------

# Row 16 - Rape
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 16.666666666666
$ws.Range("I16").Value = 34
$ws.Range("J16").Value = 50
$ws.Range("K16").Value = -32
$ws.Range("L16").Value = -27.659574468085

# Row 17 - Robbery
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 145
$ws.Range("J17").Value = 136
$ws.Range("K17").Value = 6.617647058823
$ws.Range("L17").Value = 2.836879432624

# Row 18 - Fel. Assault
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = -33.333333333333
$ws.Range("I18").Value = 27
$ws.Range("K18").Value = -41.304347826087
$ws.Range("L18").Value = -37.209302325581

# Row 19 - Burglary
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -100
$ws.Range("F19").Value = 12
$ws.Range("H19").Value = -40
$ws.Range("J19").Value = 151
$ws.Range("K19").Value = -19.205298013245
$ws.Range("L19").Value = -22.292993630573

# Row 20 - Gr. Larceny
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = -50
$ws.Range("I20").Value = 103
$ws.Range("J20").Value = 97
$ws.Range("K20").Value = 6.185567010309
$ws.Range("L20").Value = -0.961538461538

# Row 21 - G.L.A.
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = -33.333333333333
$ws.Range("F21").Value = 55
$ws.Range("G21").Value = 72
$ws.Range("H21").Value = -23.611111111111
$ws.Range("I21").Value = 451
$ws.Range("J21").Value = 487
$ws.Range("K21").Value = -7.392197125256
$ws.Range("L21").Value = -9.619238476953

# Row 24 - Transit
$ws.Range("C24").Value = 10
$ws.Range("E24").Value = -9.090909090909
$ws.Range("F24").Value = 64
$ws.Range("G24").Value = 44
$ws.Range("H24").Value = 45.454545454545
$ws.Range("I24").Value = 301
$ws.Range("J24").Value = 336
$ws.Range("K24").Value = -10.416666666666
$ws.Range("L24").Value = -7.384615384615

# Row 25 - Housing
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 3
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = -75
$ws.Range("I25").Value = 60
$ws.Range("J25").Value = 91
$ws.Range("K25").Value = -34.065934065934
$ws.Range("L25").Value = -7.692307692307

# Row 26 - Petit Larceny
$ws.Range("C26").Value = 15
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 54
$ws.Range("G26").Value = 48
$ws.Range("H26").Value = 12.5
$ws.Range("I26").Value = 271
$ws.Range("J26").Value = 263
$ws.Range("K26").Value = 3.041825095057
$ws.Range("L26").Value = 29.665071770334

# Row 28 - Misd. Assault
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 4
$ws.Range("I28").Value = 12
$ws.Range("K28").Value = 33.333333333333
$ws.Range("L28").Value = -33.333333333333

# Row 29 - Shooting Vic.
$ws.Range("G29").Value = 2
$ws.Range("J29").Value = 8
$ws.Range("K29").Value = -87.5

# Row 30 - Shooting Inc.
$ws.Range("G30").Value = 2
$ws.Range("J30").Value = 6
$ws.Range("K30").Value = -83.333333333333
